$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.693.57'
$ws.Range('E2').Value = '  +0.28%  '

$ws.Range('D3').Value = '2.525.70'
$ws.Range('E3').Value = '  +0.51%  '

$ws.Range('D5').Value = '''316.65'
$ws.Range('E5').Value = '  -0.42%  '

$ws.Range('D6').Value = '''96.88'
$ws.Range('E6').Value = '  +1.91%  '

$ws.Range('D7').Value = '''0.576'
$ws.Range('E7').Value = '  -0.24%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('D9').Value = '''0.532'
$ws.Range('E9').Value = '  -0.61%  '

$ws.Range('D10').Value = '''35.72'
$ws.Range('E10').Value = '  -0.33%  '

$ws.Range('D11').Value = '''0.0807'
$ws.Range('E11').Value = '  +0.07%  '

$ws.Range('D12').Value = '''7.52'
$ws.Range('E12').Value = '  -0.16%  '

$ws.Range('E13').Value = '  -2.51%  '

$ws.Range('D14').Value = '2.911.66'
$ws.Range('E14').Value = '  +0.60%  '

$ws.Range('D15').Value = '2.536.94'
$ws.Range('E15').Value = '  +1.09%  '

$ws.Range('D16').Value = '''15.11'
$ws.Range('E16').Value = '  -2.46%  '

$ws.Range('D17').Value = '''0.849'
$ws.Range('E17').Value = '  -0.66%  '

$ws.Range('D18').Value = '42.736.71'
$ws.Range('E18').Value = '  +0.38%  '

$ws.Range('D19').Value = '''6.82'
$ws.Range('E19').Value = '  +4.44%  '

$ws.Range('D20').Value = '''12.75'
$ws.Range('E20').Value = '  -3.07%  '

$ws.Range('D21').Value = '0.0₃0961'
$ws.Range('E21').Value = '  -0.59%  '

$ws.Range('D22').Value = '''69.69'
$ws.Range('E22').Value = '  -2.16%  '

$ws.Range('D23').Value = '''251.21'
$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('E24').Value = '  -1.75%  '

$ws.Range('E25').Value = '  +0.49%  '

$ws.Range('D26').Value = '''26.46'
$ws.Range('E26').Value = '  -0.82%  '

$ws.Range('E27').Value = '  +0.02%  '

$ws.Range('D28').Value = '''2.41'
$ws.Range('E28').Value = '  +2.15%  '

$ws.Range('D29').Value = '''40.74'
$ws.Range('E29').Value = '  +4.53%  '

$ws.Range('D30').Value = '''10.37'
$ws.Range('E30').Value = '  +3.28%  '

$ws.Range('D31').Value = '''5.92'
$ws.Range('E31').Value = '  +0.90%  '

$ws.Range('D32').Value = '''158.15'
$ws.Range('E32').Value = '  +1.48%  '

$ws.Range('D33').Value = '''2.15'
$ws.Range('E33').Value = '  +3.48%  '

$ws.Range('E34').Value = '  +4.02%  '

$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('D36').Value = '''18.90'
$ws.Range('E36').Value = '  -3.08%  '

$ws.Range('D37').Value = '''0.0785'
$ws.Range('E37').Value = '  +0.21%  '

$ws.Range('E38').Value = '  -0.21%  '

$ws.Range('E39').Value = '  -0.79%  '

$ws.Range('D40').Value = '''2.32'
$ws.Range('E40').Value = '  +10.07%  '

$ws.Range('D41').Value = '''22.38'
$ws.Range('E41').Value = '  -6.51%  '

$ws.Range('E42').Value = '  -0.69%  '

$ws.Range('D43').Value = '''0.0305'
$ws.Range('E43').Value = '  +1.61%  '

$ws.Range('E44').Value = '  +0.09%  '

$ws.Range('D45').Value = '2.029.87'
$ws.Range('E45').Value = '  -0.71%  '

$ws.Range('D46').Value = '''3.26'
$ws.Range('E46').Value = '  -2.92%  '

$ws.Range('D47').Value = '''9.06'
$ws.Range('E47').Value = '  +2.98%  '

$ws.Range('D48').Value = '''84.32'
$ws.Range('E48').Value = '  +0.08%  '

$ws.Range('D49').Value = '''105.88'
$ws.Range('E49').Value = '  +4.59%  '

$ws.Range('D50').Value = '''75.26'
$ws.Range('E50').Value = '  +3.78%  '

$ws.Range('D51').Value = '2.766.72'
$ws.Range('E51').Value = '  +0.47%  '
